$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 112144586
$ws.Range("B3").Value = 44322
$ws.Range("C3").Value = "Ovaliderad"
$ws.Range("D3").Value = "NT"
$ws.Range("E3").Value = 102366
$ws.Range("F3").Value = "Ängsmetallvinge"
$ws.Range("G3").Value = "Adscita statices"
$ws.Range("H3").Value = "(Linnaeus, 1758)"

# I3 is an empty (blank) string cell, like I2 in the template.
$ws.Range("I3").Value = "x"
$ws.Range("I3").Value = ""
$ws.Range("I3").Style = "Normal"

$ws.Range("P3").Value = "Vä, delomr 10, 450 m ONO fängelset, Sk"
$ws.Range("Q3").Value = 442994.5965538067
$ws.Range("R3").Value = 6204827.083255709
$ws.Range("S3").Value = 10
$ws.Range("T3").Value = "Skåne"
$ws.Range("U3").Value = "Kristianstad"
$ws.Range("V3").Value = "Skåne"
$ws.Range("W3").Value = "Vä"

# Date-looking text values must stay text, not get converted to Excel date serials.
$ws.Range("Y3").NumberFormat = "@"
$ws.Range("Y3").Value = "2013-06-03"
$ws.Range("Y3").Style = "Normal"

$ws.Range("Z3").Value = "00:00"

$ws.Range("AA3").NumberFormat = "@"
$ws.Range("AA3").Value = "2013-06-03"
$ws.Range("AA3").Style = "Normal"

$ws.Range("AB3").Value = "00:00"

$ws.Range("AD3").Value = $false
$ws.Range("AE3").Value = $false
$ws.Range("AG3").Value = $false
$ws.Range("AI3").Value = "på igenväxande grässandmark"

# AT3 is an empty (blank) string cell, like AT2 in the template.
$ws.Range("AT3").Value = "x"
$ws.Range("AT3").Value = ""
$ws.Range("AT3").Style = "Normal"

$ws.Range("AW3").Value = "Nils Otto Nilsson"
$ws.Range("AX3").Value = "Nils Otto Nilsson"
$ws.Range("AY3").Value = "Krst NV-program 2013"
